# Fix bugs: - Read formula
#
# Renames Sheet2 -> "Data", fills it with sample rows covering an Int,
# Double, Double Formula, String, Percent, String Formular, Date,
# Date Func and Date Formular, and leaves that sheet active/selected
# (matching the "Read formula" sample workbook used by the unit tests).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "Data"

# Row 1 - plain integer
$ws.Range("A1").Value = "Int"
$ws.Range("B1").Value = 1

# Row 2 - plain double
$ws.Range("A2").Value = "Double"
$ws.Range("B2").Value = 0.25

# Row 3 - formula returning a double
$ws.Range("A3").Value = "Double Formula"
$ws.Range("B3").Formula = "=6/10"

# Row 4 - plain string
$ws.Range("A4").Value = "String"
$ws.Range("B4").Value = "Test"

# Row 5 - percentage number format
$ws.Range("A5").Value = "Percent"
$ws.Range("B5").NumberFormat = "0%"
$ws.Range("B5").Value = 0.1

# Row 6 - formula returning a string
$ws.Range("A6").Value = "String Formular"
$ws.Range("B6").Formula = '="A" & "B"'

# Row 7 - literal date value
$ws.Range("A7").Value = "Date"
$ws.Range("B7").NumberFormat = "d-mmm"
$ws.Range("B7").Value = 40939

# Labels for rows 8/9 are entered out of row order so the shared-string
# table ends up with "Date Formular" before "Date Func", matching the
# original authoring order.
$ws.Range("A9").Value = "Date Formular"
$ws.Range("A8").Value = "Date Func"

# Row 8 - volatile date formula
$ws.Range("B8").NumberFormat = "mm-dd-yy"
$ws.Range("B8").Formula = "=TODAY()"

# Row 9 - formula referencing the other date formula
$ws.Range("B9").NumberFormat = "mm-dd-yy"
$ws.Range("B9").Formula = "=B8+1"

# Column widths matching the original sheet
$ws.Columns.Item(1).ColumnWidth = 14.451822916666666
$ws.Columns.Item(2).ColumnWidth = 13.022135416666666

# Leave "Data" as the active sheet/tab with the same selection it was
# saved with.
$ws.Range("G19").Select() | Out-Null
